$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column U (21st column). Excel's default
# column-insert behavior shifts everything from U onward one column to the
# right and formats the freshly inserted column like the column to its
# left (matches what happens when you right-click a column header and
# choose "Insert").
$ws.Columns("U:U").Insert()

# Give the new column its header text.
$ws.Range("U1").Value = "Sub brand"

# The AutoFilter / _FilterDatabase range used to stop at column AO; now
# that there is one more column it must cover through AP. Re-apply the
# AutoFilter over the corrected range (toggle off first since the range
# object still references the stale AO boundary).
$ws.AutoFilterMode = $false
$ws.Range("A1:AP37").AutoFilter()

# Update the defined names backing the filter ("_FilterDatabase" and
# "_FilterDatabase_0") so they also point at the widened range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='HoReCa Bar Tavern_Night Club'!`$A`$1:`$AP`$37"
    }
}

# Restore the view focus/selection roughly where the author left it.
$ws.Range("T31").Select()
